$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new backlog rows (10 paintings / image work added in Firebase) ---
# Row 41: "Add Image Reference" task with comment about image source/credit
$ws.Range("A41").Value = "Add Image Reference"
$ws.Range("E41").Value = "Images and Descriptions from National Gallery of Art http://www.nga.gov/"

# Row 42: "Add Images in Firebase" task, owner DM, status Done
$ws.Range("A42").Value = "Add Images in Firebase"
$ws.Range("C42").Value = "DM"
$ws.Range("D42").Value = "Done"

# --- Highlight the "Later" section header (row 44) ---
# Bold + yellow fill for the section label itself
$ws.Range("A44").Font.Bold = $true
$ws.Range("A44").Interior.Color = 65535

# Extend the yellow highlight band across the rest of the row (B:G)
$ws.Range("B44:G44").Interior.Color = 65535

# --- Highlight the table header row to match ---
$ws.Range("A1:E1").Interior.Color = 65535

# --- Restore the selection to where the author left off editing ---
[void]$ws.Range("H17").Select()
